# Refresh the cryptos list (Coin / Link / Price / Volume(1h)) with the
# latest scrape -- see commit message "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> new text value.
$updates = [ordered]@{
    "D2" = "49.636.61"
    "E2" = "  +2.98%  "
    "D3" = "2.553.67"
    "E3" = "  +1.65%  "
    "E4" = "  -0.09%  "
    "D5" = "322.41"
    "E5" = "  +0.19%  "
    "D6" = "108.74"
    "E6" = "  +0.11%  "
    "D7" = "0.529"
    "E7" = "  +0.17%  "
    "E8" = "  -0.03%  "
    "D9" = "0.556"
    "E9" = "  +2.18%  "
    "E10" = "  +0.27%  "
    "D11" = "20.46"
    "E11" = "  +1.25%  "
    "D12" = "0.0817"
    "E12" = "  -0.38%  "
    "E13" = "  +0.47%  "
    "D14" = "7.28"
    "E14" = "  +0.90%  "
    "D15" = "2.950.29"
    "E15" = "  +1.68%  "
    "D16" = "2.549.34"
    "E16" = "  +1.64%  "
    "D17" = "0.863"
    "E17" = "  +2.02%  "
    "D18" = "49.421.28"
    "E18" = "  +2.90%  "
    "E19" = "  +10.52%  "
    "D20" = "13.24"
    "E20" = "  +0.75%  "
    "D21" = "6.67"
    "E21" = "  -0.02%  "
    "E22" = "  -1.15%  "
    "D23" = "283.41"
    "E23" = "  +2.14%  "
    "D24" = "72.17"
    "E24" = "  -0.05%  "
    "E25" = "  -1.30%  "
    "D26" = "26.43"
    "E26" = "  +1.99%  "
    "E27" = "  -0.07%  "
    "E28" = "  -1.27%  "
    "D29" = "0.144"
    "E29" = "  +4.09%  "
    "D30" = "9.85"
    "E30" = "  -0.21%  "
    "D31" = "35.27"
    "E31" = "  -0.29%  "
    "D32" = "49.62"
    "E32" = "  +0.73%  "
    "D33" = "19.75"
    "E33" = "  +1.57%  "
    "D34" = "5.40"
    "E34" = "  +0.77%  "
    "E35" = "  -0.26%  "
    "D36" = "0.0787"
    "E36" = "  -0.01%  "
    "D37" = "2.04"
    "E37" = "  +3.98%  "
    "D38" = "4.71"
    "E38" = "  +2.10%  "
    "D39" = "3.01"
    "E39" = "  +1.84%  "
    "E40" = "  +0.18%  "
    "D41" = "121.78"
    "E41" = "  -1.04%  "
    "B42" = "EnergySwap"
    "C42" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D42" = "22.29"
    "E42" = "  +3.95%  "
    "B43" = "WEMIXToken"
    "C43" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D43" = "2.21"
    "E43" = "  -0.31%  "
    "D44" = "0.0312"
    "E44" = "  +1.92%  "
    "E45" = "  +6.16%  "
    "D46" = "2.017.24"
    "E46" = "  +0.71%  "
    "D47" = "2.01"
    "E47" = "  +7.98%  "
    "D48" = "2.14"
    "E48" = "  +7.99%  "
    "D49" = "9.01"
    "E49" = "  +0.07%  "
    "D50" = "5.32"
    "E50" = "  +1.94%  "
    "D51" = "81.37"
    "E51" = "  +1.23%  "
}

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "20.46", "5.40"), silently dropping trailing zeros / exact formatting.
# These are written with the cell temporarily forced to Text format, then the
# original style is restored so no stray formatting is left behind.
$numericLooking = @(
    "D5"
    "D6"
    "D7"
    "D9"
    "D11"
    "D12"
    "D14"
    "D17"
    "D20"
    "D21"
    "D23"
    "D24"
    "D26"
    "D29"
    "D30"
    "D31"
    "D32"
    "D33"
    "D34"
    "D36"
    "D37"
    "D38"
    "D39"
    "D41"
    "D42"
    "D43"
    "D44"
    "D47"
    "D48"
    "D49"
    "D50"
    "D51"
)

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $range = $ws.Range($addr)

    if ($numericLooking -contains $addr) {
        $originalStyle = $range.Style
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = $originalStyle
    } else {
        $range.Value = $value
    }
}
